$wb = $excel.ActiveWorkbook

# --- Sheet "설치" (install) gets two new rows with content ---
$ws2 = $wb.Worksheets.Item("설치")

$ws2.Range("A6").Value = "환경변수"
$ws2.Range("B7").Value = "https://mongodev.tistory.com/28"
$ws2.Range("B7").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://mongodev.tistory.com/28") | Out-Null

# page setup triggers printer settings relationship like sheet1
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- view state ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$ws2.Activate()
$ws2.Range("G10").Select()
